$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data rows (rows 2-6), columns A-J
# A=trialTrain, B=x_fixStart, C=y_fixStart, D=x_corrSteps, E=y_corrSteps,
# F=x_nrSteps, G=y_nrSteps, H=alienID, I=praclen, J=version
$data = @(
    @(1, 8, 7, 4, 5, -4, -2, 23, 5, "train_dim2_1"),
    @(2, 7, 7, 2, 6, -5, -1, 12, 5, "train_dim2_1"),
    @(3, 6, 8, 5, 3, -1, -5, 56, 5, "train_dim2_1"),
    @(4, 9, 9, 7, 5, -2, -4, 45, 5, "train_dim2_1"),
    @(5, 7, 5, 4, 2, -3, -3, 34, 5, "train_dim2_1")
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le 10; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
